$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (2024-09-21)
$ws.Range("B2").Value = 3.230985683306322
$ws.Range("C2").Value = 1.667794583268128
$ws.Range("D2").Value = 0.8054896365839992
$ws.Range("E2").Value = 8.660232485948974
$ws.Range("G2").Value = 14.36450238910742

# Row 3 (2024-09-17)
$ws.Range("B3").Value = 0.3048080303191223
$ws.Range("C3").Value = 1.667794583268128
$ws.Range("D3").Value = 0.8054896365839992
$ws.Range("E3").Value = 0.496779210170732
$ws.Range("G3").Value = 3.274871460341982
